$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column A. This shifts the existing
#    PRODUCT_ID..CAPACITY_SHIFT_3 columns from A:F to B:G, each cell keeping
#    its original style (header style s=1, data style s=2).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).Insert()

# ---------------------------------------------------------------------------
# 2. Extend the table with 4 extra data rows (8-11), matching the existing
#    blank-row formatting, by copying row 7's formats down.
# ---------------------------------------------------------------------------
$ws.Range("A7:G7").Copy()
$ws.Range("A8:G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. New column A data cells (rows 2-11): give them the same "boxed" border
#    style used by the rest of the data cells, by copying the format from an
#    existing data cell (C2, style s=2) onto A2:A11.
# ---------------------------------------------------------------------------
$ws.Range("C2").Copy()
$ws.Range("A2:A11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Column B header (B1, "PRODUCT_ID"): same header look as before, but the
#    left border is removed now that column A sits to its left.
# ---------------------------------------------------------------------------
$ws.Range("B1").Borders.Item(7).LineStyle = -4142

# ---------------------------------------------------------------------------
# 5. Column B data cells (rows 2-11): drop the left border edge so column A's
#    right edge is the only divider between the two columns.
# ---------------------------------------------------------------------------
$ws.Range("B2:B11").Borders.Item(7).LineStyle = -4142

# ---------------------------------------------------------------------------
# 6. Column A header (A1): "NOMOR" label with its own bold black Aptos Narrow
#    font, the same orange header fill, and a thin black box border.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "NOMOR"
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Color = 0
$ws.Range("A1").Interior.Color = 49407
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").Borders.Color = 0
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").Borders.Weight = 2

# ---------------------------------------------------------------------------
# 7. Autofit the columns that should hug their contents.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).EntireColumn.AutoFit()
$ws.Columns.Item(2).EntireColumn.AutoFit()
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(7).EntireColumn.AutoFit()

$ws.Range("A1").Select()
